# Insert one new data row at row 515 (pushes existing rows 515:606 down to 516:607)
# and populate it with a new weekly price record for the Brócoli sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 515..606 down by one row.
$ws.Rows.Item(515).EntireRow.Insert()

# Populate the newly inserted row 515 with the new record.
$ws.Range("A515").Value = 4
$ws.Range("B515").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C515").Value = "Los Lagos"
$ws.Range("D515").Value = 45209
$ws.Range("E515").Value = 10
$ws.Range("F515").Value = 100112023
$ws.Range("G515").Value = "Brócoli"
$ws.Range("H515").Value = "Sin especificar"
$ws.Range("I515").Value = "Primera"
$ws.Range("J515").Value = 1500
$ws.Range("K515").Value = 1500
$ws.Range("L515").Value = 1500
$ws.Range("M515").Value = 1500
$ws.Range("N515").Value = "$/unidad"
$ws.Range("O515").Value = "Región Metropolitana"
$ws.Range("P515").Value = 1500
$ws.Range("Q515").Value = 1
$ws.Range("R515").Value = "Hortaliza"
